$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.1249546554732959
$ws.Range("D2").Value = 0.07506655530873729
$ws.Range("E2").Value = 0.1169651680031922
$ws.Range("F2").Value = 2.544638687492935
$ws.Range("G2").Value = 1.862070362332105
$ws.Range("H2").Value = 1.615786594251816
$ws.Range("I2").Value = 1.81938182001209
$ws.Range("J2").Value = 0.1838533104454285
$ws.Range("K2").Value = 1.848937505168578
$ws.Range("L2").Value = 0.1811597084840884
$ws.Range("M2").Value = 0.5601606953080207
$ws.Range("N2").Value = 1.898674655541754
$ws.Range("C3").Value = 0.1238174759113591
$ws.Range("D3").Value = 0.07415474852088977
$ws.Range("E3").Value = 0.1171315795202101
$ws.Range("F3").Value = 2.550860119042454
$ws.Range("G3").Value = 1.866020650484188
$ws.Range("H3").Value = 1.623784785589635
$ws.Range("I3").Value = 1.824092095109528
$ws.Range("J3").Value = 0.1849019369038238
$ws.Range("K3").Value = 1.743109116010999
$ws.Range("L3").Value = 0.1819805800203511
$ws.Range("M3").Value = 0.5411742080338655
$ws.Range("N3").Value = 1.917608046639511
$ws.Range("C4").Value = 0.1231576442727516
$ws.Range("D4").Value = 0.07361085567465864
$ws.Range("E4").Value = 0.1172664538558728
$ws.Range("F4").Value = 2.555955896324335
$ws.Range("G4").Value = 1.869514367337644
$ws.Range("H4").Value = 1.629411293329952
$ws.Range("I4").Value = 1.827907525011085
$ws.Range("J4").Value = 0.1856015611118771
$ws.Range("K4").Value = 1.678759324706846
$ws.Range("L4").Value = 0.1825309249918128
$ws.Range("M4").Value = 0.5297281085374266
$ws.Range("N4").Value = 1.929883204265575
$ws.Range("C5").Value = 0.1228984578522727
$ws.Range("D5").Value = 0.07339326069487129
$ws.Range("E5").Value = 0.1173296589518227
$ws.Range("F5").Value = 2.558353223306426
$ws.Range("G5").Value = 1.871206451284905
$ws.Range("H5").Value = 1.6318841044626
$ws.Range("I5").Value = 1.829694457451097
$ws.Range("J5").Value = 0.1859007030925586
$ws.Range("K5").Value = 1.652695740308644
$ws.Range("L5").Value = 0.1827668659967419
$ws.Range("M5").Value = 0.5251172456733784
$ws.Range("N5").Value = 1.935048938618984
$ws.Range("C6").Value = 0.1228560073209692
$ws.Range("D6").Value = 0.07335737452680036
$ws.Range("E6").Value = 0.1173406525620937
$ws.Range("F6").Value = 2.558770669087096
$ws.Range("G6").Value = 1.871503621999679
$ws.Range("H6").Value = 1.632305582967021
$ws.Range("I6").Value = 1.830005193691697
$ws.Range("J6").Value = 0.1859512238489618
$ws.Range("K6").Value = 1.648377570894695
$ws.Range("L6").Value = 0.1828067493607222
$ws.Range("M6").Value = 0.5243548561074647
$ws.Range("N6").Value = 1.935916582376741
$ws.Range("C7").Value = 0.1231541094585751
$ws.Range("D7").Value = 0.07360790468835177
$ws.Range("E7").Value = 0.1172672728596158
$ws.Range("F7").Value = 2.555986928858601
$ws.Range("G7").Value = 1.869536101102369
$ws.Range("H7").Value = 1.62944391383374
$ws.Range("I7").Value = 1.827930684467802
$ws.Range("J7").Value = 0.1856055385772066
$ws.Range("K7").Value = 1.678407174720689
$ws.Range("L7").Value = 0.1825340596901146
$ws.Range("M7").Value = 0.5296657077646429
$ws.Range("N7").Value = 1.929952209043048
$ws.Range("C8").Value = 0.1245546115778922
$ws.Range("D8").Value = 0.07474887081955472
$ws.Range("E8").Value = 0.1170157727416861
$ws.Range("F8").Value = 2.546519009115613
$ws.Range("G8").Value = 1.863210558068673
$ws.Range("H8").Value = 1.618395862972775
$ws.Range("I8").Value = 1.820814252615399
$ws.Range("J8").Value = 0.1842033135150487
$ws.Range("K8").Value = 1.812317979938598
$ws.Range("L8").Value = 0.1814331429406124
$ws.Range("M8").Value = 0.5535703848065339
$ws.Range("N8").Value = 1.90506802535916
$ws.Range("C9").Value = 0.1276039899410364
$ws.Range("D9").Value = 0.07711166833826155
$ws.Range("E9").Value = 0.1167811643025569
$ws.Range("F9").Value = 2.53808016458332
$ws.Range("G9").Value = 1.859295325353727
$ws.Range("H9").Value = 1.602408361724443
$ws.Range("I9").Value = 1.814189334050113
$ws.Range("J9").Value = 0.1818952620481866
$ws.Range("K9").Value = 2.079869853577634
$ws.Range("L9").Value = 0.1796408846926347
$ws.Range("M9").Value = 0.6021172815760991
$ws.Range("N9").Value = 1.86142417901382
$ws.Range("C10").Value = 0.1300271492730758
$ws.Range("D10").Value = 0.07892256231782113
$ws.Range("E10").Value = 0.1167653771946942
$ws.Range("F10").Value = 2.538064662788756
$ws.Range("G10").Value = 1.861615389296446
$ws.Range("H10").Value = 1.594124715054505
$ws.Range("I10").Value = 1.813799467785088
$ws.Range("J10").Value = 0.1804678218740854
$ws.Range("K10").Value = 2.279428344497489
$ws.Range("L10").Value = 0.1785463734130488
$ws.Range("M10").Value = 0.6387939745794142
$ws.Range("N10").Value = 1.832495894331331
$ws.Range("C11").Value = 0.1311688117695695
$ws.Range("D11").Value = 0.07976236228527966
$ws.Range("E11").Value = 0.1167919784544349
$ws.Range("F11").Value = 2.539403072929545
$ws.Range("G11").Value = 1.863804090724784
$ws.Range("H11").Value = 1.591108469865759
$ws.Range("I11").Value = 1.81459643590243
$ws.Range("J11").Value = 0.1798764932557404
$ws.Range("K11").Value = 2.370856334994528
$ws.Range("L11").Value = 0.1780964532944154
$ws.Range("M11").Value = 0.6556967564847511
$ws.Range("N11").Value = 1.820015801395158
$ws.Range("C12").Value = 0.1316067500913647
$ws.Range("D12").Value = 0.08008264655316566
$ws.Range("E12").Value = 0.1168068904544999
$ws.Range("F12").Value = 2.540103526277179
$ws.Range("G12").Value = 1.864796210507393
$ws.Range("H12").Value = 1.590074457689951
$ws.Range("I12").Value = 1.815038467188714
$ws.Range("J12").Value = 0.1796609000103011
$ws.Range("K12").Value = 2.405569983361943
$ws.Range("L12").Value = 0.1779329586916916
$ws.Range("M12").Value = 0.6621285807783863
$ws.Range("N12").Value = 1.815387593976851
$ws.Range("C13").Value = 0.131512183180206
$ws.Range("D13").Value = 0.08001356715456609
$ws.Range("E13").Value = 0.1168034639921984
$ws.Range("F13").Value = 2.539944057088505
$ws.Range("G13").Value = 1.864575271236561
$ws.Range("H13").Value = 1.590292339061563
$ws.Range("I13").Value = 1.81493702866544
$ws.Range("J13").Value = 0.1797069616009459
$ws.Range("K13").Value = 2.39808971358832
$ws.Range("L13").Value = 0.1779678644717038
$ws.Range("M13").Value = 0.6607419934724916
$ws.Range("N13").Value = 1.816380014208214
$ws.Range("C14").Value = 0.1312047289250131
$ws.Range("D14").Value = 0.07978866695233933
$ws.Range("E14").Value = 0.1167931084091283
$ws.Range("F14").Value = 2.539456818210013
$ws.Range("G14").Value = 1.863882438163358
$ws.Range("H14").Value = 1.591021232883861
$ws.Range("I14").Value = 1.814629990771508
$ws.Range("J14").Value = 0.1798585893979165
$ws.Range("K14").Value = 2.373710414316861
$ws.Range("L14").Value = 0.1780828647092854
$ws.Range("M14").Value = 0.6562252847727592
$ws.Range("N14").Value = 1.819633077149682
$ws.Range("C15").Value = 0.1310171343290989
$ws.Range("D15").Value = 0.07965120377048862
$ws.Range("E15").Value = 0.1167873948852201
$ws.Range("F15").Value = 2.539183590741189
$ws.Range("G15").Value = 1.863479335304618
$ws.Range("H15").Value = 1.591481789940531
$ws.Range("I15").Value = 1.81446018808272
$ws.Range("J15").Value = 0.1799525502319668
$ws.Range("K15").Value = 2.358789312924785
$ws.Range("L15").Value = 0.1781542011471622
$ws.Range("M15").Value = 0.6534627127237087
$ws.Range("N15").Value = 1.821638400470007
$ws.Range("C16").Value = 0.1299533271274527
$ws.Range("D16").Value = 0.07886799887837981
$ws.Range("E16").Value = 0.1167643168885402
$ws.Range("F16").Value = 2.538004281474358
$ws.Range("G16").Value = 1.86149518625399
$ws.Range("H16").Value = 1.594336969318846
$ws.Range("I16").Value = 1.813767000972419
$ws.Range("J16").Value = 0.1805076329816941
$ws.Range("K16").Value = 2.273466234587318
$ws.Range("L16").Value = 0.1785767405847878
$ws.Range("M16").Value = 0.6376937070415707
$ws.Range("N16").Value = 1.833325172963605
$ws.Range("C17").Value = 0.129310764479186
$ws.Range("D17").Value = 0.078391606567358
$ws.Range("E17").Value = 0.1167587981483642
$ws.Range("F17").Value = 2.537625524655894
$ws.Range("G17").Value = 1.860568479547041
$ws.Range("H17").Value = 1.596281164942823
$ws.Range("I17").Value = 1.813591392118724
$ws.Range("J17").Value = 0.1808630093395287
$ws.Range("K17").Value = 2.221288283684146
$ws.Range("L17").Value = 0.1788482298166869
$ws.Range("M17").Value = 0.6280756530341733
$ws.Range("N17").Value = 1.840668683440107
$ws.Range("C18").Value = 0.1289448863314249
$ws.Range("D18").Value = 0.0781191076804717
$ws.Range("E18").Value = 0.1167588048639594
$ws.Range("F18").Value = 2.537534305751223
$ws.Range("G18").Value = 1.860142112720354
$ws.Range("H18").Value = 1.597470197017657
$ws.Range("I18").Value = 1.813582086192923
$ws.Range("J18").Value = 0.1810728744210515
$ws.Range("K18").Value = 2.191337960820761
$ws.Range("L18").Value = 0.1790089003283768
$ws.Range("M18").Value = 0.6225641742350092
$ws.Range("N18").Value = 1.844956435438753
$ws.Range("C19").Value = 0.1288216442755612
$ws.Range("D19").Value = 0.07802710447565886
$ws.Range("E19").Value = 0.1167593540578231
$ws.Range("F19").Value = 2.537525165241306
$ws.Range("G19").Value = 1.860016059236074
$ws.Range("H19").Value = 1.597884938769838
$ws.Range("I19").Value = 1.813594680934976
$ws.Range("J19").Value = 0.1811448696437381
$ws.Range("K19").Value = 2.181207822187332
$ws.Range("L19").Value = 0.1790640769946847
$ws.Range("M19").Value = 0.6207016243698504
$ws.Range("N19").Value = 1.846419179190221
$ws.Range("C20").Value = 0.1293787830379927
$ws.Range("D20").Value = 0.07844216336653176
$ws.Range("E20").Value = 0.1167590565514161
$ws.Range("F20").Value = 2.537652736940942
$ws.Range("G20").Value = 1.860656088686241
$ws.Range("H20").Value = 1.596066876066729
$ws.Range("I20").Value = 1.813600594467928
$ws.Range("J20").Value = 0.1808246137197678
$ws.Range("K20").Value = 2.226836406206644
$ws.Range("L20").Value = 0.1788188619601492
$ws.Range("M20").Value = 0.6290973843657213
$ws.Range("N20").Value = 1.83988033458423
$ws.Range("C21").Value = 0.1312948836358885
$ws.Range("D21").Value = 0.07985466423925658
$ws.Range("E21").Value = 0.1167960189273209
$ws.Range("F21").Value = 2.539594675894051
$ws.Range("G21").Value = 1.864081505228739
$ws.Range("H21").Value = 1.590804202995599
$ws.Range("I21").Value = 1.814716368203946
$ws.Range("J21").Value = 0.1798138266679565
$ws.Range("K21").Value = 2.380868722976459
$ws.Range("L21").Value = 0.1780488997698342
$ws.Range("M21").Value = 0.6575511096836379
$ws.Range("N21").Value = 1.818674921093759
$ws.Range("C22").Value = 0.1325798775447993
$ws.Range("D22").Value = 0.08079103293007961
$ws.Range("E22").Value = 0.1168483706684746
$ws.Range("F22").Value = 2.54199252095664
$ws.Range("G22").Value = 1.867272245344736
$ws.Range("H22").Value = 1.587995282545421
$ws.Range("I22").Value = 1.816263060454105
$ws.Range("J22").Value = 0.1792017668902623
$ws.Range("K22").Value = 2.482072534531142
$ws.Range("L22").Value = 0.1775857804605074
$ws.Range("M22").Value = 0.6763284481839094
$ws.Range("N22").Value = 1.805385632412559
$ws.Range("C23").Value = 0.1318910733901504
$ws.Range("D23").Value = 0.08029007666458199
$ws.Range("E23").Value = 0.1168178560703659
$ws.Range("F23").Value = 2.540609421334295
$ws.Range("G23").Value = 1.8654820632363
$ws.Range("H23").Value = 1.58943674997866
$ws.Range("I23").Value = 1.815362719330167
$ws.Range("J23").Value = 0.1795239969076086
$ws.Range("K23").Value = 2.428009651270713
$ws.Range("L23").Value = 0.1778292934036152
$ws.Range("M23").Value = 0.6662901474127452
$ws.Range("N23").Value = 1.812426243181925
$ws.Range("C24").Value = 0.1293480208196058
$ws.Range("D24").Value = 0.07841930232053329
$ws.Range("E24").Value = 0.1167589298222893
$ws.Range("F24").Value = 2.537640040123108
$ws.Range("G24").Value = 1.860616149153515
$ws.Range("H24").Value = 1.596163534001306
$ws.Range("I24").Value = 1.813596148591124
$ws.Range("J24").Value = 0.1808419550765912
$ws.Range("K24").Value = 2.224327952029682
$ws.Range("L24").Value = 0.1788321248760418
$ws.Range("M24").Value = 0.6286354033260508
$ws.Range("N24").Value = 1.84023654191175
$ws.Range("C25").Value = 0.1267468152919236
$ws.Range("D25").Value = 0.07645920157386854
$ws.Range("E25").Value = 0.1168170590640063
$ws.Range("F25").Value = 2.53927779342699
$ws.Range("G25").Value = 1.859443479737479
$ws.Range("H25").Value = 1.606125463860664
$ws.Range("I25").Value = 1.815195896490032
$ws.Range("J25").Value = 0.182472466538588
$ws.Range("K25").Value = 2.00696299749734
$ws.Range("L25").Value = 0.1800866170159345
$ws.Range("M25").Value = 0.5888061954995294
$ws.Range("N25").Value = 1.872679785247346

Write-Host "Updated pl_mw values for 380 kV case"
